# Apply the "456a3b4" data refresh to 江西-漫展信息.xlsx
#
# Sheets in the workbook:
#   1 = 展览     (Exhibition)   -> last data row moves from 36 to 37
#   2 = 演出     (Performance)  -> unchanged
#   3 = 本地生活 (Local life)   -> unchanged
#   4 = 全部类型 (All types)    -> last data row moves from 37 to 38
#
# For both sheet 1 and sheet 4 we:
#   - bump a handful of "want to go" counts (column F)
#   - turn G14 ("min ticket price") from a number into the text "已售罄"
#     (sold out)
#   - append one brand-new row for "南昌·萌卡动漫展" on 2025-01-01
#
# NOTE 1: this COM-interop host loses COM object values when they are passed
# into functions via named parameters (-ws $ws); positional parameters work
# fine, so every helper below uses positional parameters only.
#
# NOTE 2: this host also mishandles parenthesized expressions used directly
# as positional call arguments (e.g. `Foo $ws ("B" + $n) "x"`); the call
# silently turns into a no-op. Always assign the expression to a local
# variable first and pass the variable.
#
# NOTE 3: assigning a date-shaped literal (e.g. "2025-01-01") straight to
# Range.Value makes Excel auto-convert it to a real date serial, which does
# not match the source data (plain text). To keep it as text we stage the
# value through a text-formula helper cell and paste-special the computed
# value (paste-values never re-triggers Excel's "smart" literal parsing).

function Update-CommonCells($ws) {
    $ws.Range("F7").Value  = 2188
    $ws.Range("F11").Value = 4958
    $ws.Range("F17").Value = 188
    $ws.Range("F21").Value = 3912
    $ws.Range("F23").Value = 679
    $ws.Range("G14").Value = "已售罄"
}

function Set-TextNoAutoConvert($ws, [string]$cellAddr, [string]$text) {
    $helper = $ws.Range("ZZ1")
    $formulaText = '="' + $text + '"'
    $helper.Formula = $formulaText
    $helper.Copy()
    $destRange = $ws.Range($cellAddr)
    $destRange.PasteSpecial(-4163)  # xlPasteValues
    $helper.ClearContents()
}

function Set-NewEventRow($ws, [int]$rowNum, [int]$index) {
    # Clone formatting (including the bold/centered/bordered style used for
    # column A) from the previous last row, then fill in the new values.
    $prevRow = $rowNum - 1
    $srcAddr = "A" + $prevRow
    $dstAddr = "A" + $rowNum
    $srcRange = $ws.Range($srcAddr)
    $dstRange = $ws.Range($dstAddr)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($rowNum, 1).Value = $index

    $bAddr = "B" + $rowNum
    Set-TextNoAutoConvert $ws $bAddr "2025-01-01"

    $ws.Cells.Item($rowNum, 3).Value = "南昌·萌卡动漫展"
    $ws.Cells.Item($rowNum, 4).Value = "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
    $ws.Cells.Item($rowNum, 5).Value = "2025.01.01 09:00-01.03 17:00"
    $ws.Cells.Item($rowNum, 6).Value = 0
    $ws.Cells.Item($rowNum, 7).Value = 65
    $ws.Cells.Item($rowNum, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93031"
    $ws.Cells.Item($rowNum, 9).Value = "//i2.hdslb.com/bfs/openplatform/202409/HTlK8fN21727112669248.jpeg"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (36 data rows -> 37)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Update-CommonCells $ws1
$ws1.Range("F34").Value = 977
$ws1.Range("F35").Value = 2486

Set-NewEventRow $ws1 37 36

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (37 data rows -> 38)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Update-CommonCells $ws4
$ws4.Range("F35").Value = 977
$ws4.Range("F36").Value = 2486

Set-NewEventRow $ws4 38 37

$wb.Save()
